$d = $word.ActiveDocument

$d.Content.Find.Execute("29-12=17", $true, $false, $false, $false, $false, $true, 1, $false, "13+0=13", 2) | Out-Null
$d.Content.Find.Execute("32-31=1", $true, $false, $false, $false, $false, $true, 1, $false, "16-11=5", 2) | Out-Null
$d.Content.Find.Execute("36+36=72", $true, $false, $false, $false, $false, $true, 1, $false, "93-4=89", 2) | Out-Null
$d.Content.Find.Execute("61-3=58", $true, $false, $false, $false, $false, $true, 1, $false, "26+9=35", 2) | Out-Null
$d.Content.Find.Execute("80-9=71", $true, $false, $false, $false, $false, $true, 1, $false, "93-71=22", 2) | Out-Null
$d.Content.Find.Execute("48+0=48", $true, $false, $false, $false, $false, $true, 1, $false, "82+2=84", 2) | Out-Null
$d.Content.Find.Execute("95-62=33", $true, $false, $false, $false, $false, $true, 1, $false, "58+3=61", 2) | Out-Null
$d.Content.Find.Execute("58-14=44", $true, $false, $false, $false, $false, $true, 1, $false, "50-14=36", 2) | Out-Null
$d.Content.Find.Execute("6+65=71", $true, $false, $false, $false, $false, $true, 1, $false, "66+17=83", 2) | Out-Null
$d.Content.Find.Execute("12+5=17", $true, $false, $false, $false, $false, $true, 1, $false, "13+18=31", 2) | Out-Null
$d.Content.Find.Execute("98-14=84", $true, $false, $false, $false, $false, $true, 1, $false, "9+40=49", 2) | Out-Null
$d.Content.Find.Execute("19+35=54", $true, $false, $false, $false, $false, $true, 1, $false, "75+8=83", 2) | Out-Null
$d.Content.Find.Execute("41+54=95", $true, $false, $false, $false, $false, $true, 1, $false, "83+1=84", 2) | Out-Null
$d.Content.Find.Execute("0+13=13", $true, $false, $false, $false, $false, $true, 1, $false, "38-29=9", 2) | Out-Null
$d.Content.Find.Execute("35-17=18", $true, $false, $false, $false, $false, $true, 1, $false, "56-15=41", 2) | Out-Null
$d.Content.Find.Execute("45-22=23", $true, $false, $false, $false, $false, $true, 1, $false, "0+90=90", 2) | Out-Null
$d.Content.Find.Execute("28+63=91", $true, $false, $false, $false, $false, $true, 1, $false, "11+66=77", 2) | Out-Null
$d.Content.Find.Execute("92-50=42", $true, $false, $false, $false, $false, $true, 1, $false, "87-19=68", 2) | Out-Null
$d.Content.Find.Execute("86-45=41", $true, $false, $false, $false, $false, $true, 1, $false, "80-17=63", 2) | Out-Null
$d.Content.Find.Execute("78-60=18", $true, $false, $false, $false, $false, $true, 1, $false, "20+34=54", 2) | Out-Null
$d.Content.Find.Execute("37-6=31", $true, $false, $false, $false, $false, $true, 1, $false, "77+12=89", 2) | Out-Null
$d.Content.Find.Execute("49+9=58", $true, $false, $false, $false, $false, $true, 1, $false, "57-37=20", 2) | Out-Null
$d.Content.Find.Execute("70-55=15", $true, $false, $false, $false, $false, $true, 1, $false, "48+48=96", 2) | Out-Null
$d.Content.Find.Execute("51+13=64", $true, $false, $false, $false, $false, $true, 1, $false, "12+79=91", 2) | Out-Null
$d.Content.Find.Execute("61+38=99", $true, $false, $false, $false, $false, $true, 1, $false, "67-42=25", 2) | Out-Null
$d.Content.Find.Execute("4+70=74", $true, $false, $false, $false, $false, $true, 1, $false, "72-20=52", 2) | Out-Null
$d.Content.Find.Execute("99-35=64", $true, $false, $false, $false, $false, $true, 1, $false, "27-22=5", 2) | Out-Null
$d.Content.Find.Execute("93-64=29", $true, $false, $false, $false, $false, $true, 1, $false, "88-40=48", 2) | Out-Null
$d.Content.Find.Execute("94-23=71", $true, $false, $false, $false, $false, $true, 1, $false, "8+70=78", 2) | Out-Null
$d.Content.Find.Execute("32-5=27", $true, $false, $false, $false, $false, $true, 1, $false, "46-24=22", 2) | Out-Null
$d.Content.Find.Execute("3+67=70", $true, $false, $false, $false, $false, $true, 1, $false, "96-67=29", 2) | Out-Null
$d.Content.Find.Execute("85-33=52", $true, $false, $false, $false, $false, $true, 1, $false, "91-46=45", 2) | Out-Null
$d.Content.Find.Execute("38+50=88", $true, $false, $false, $false, $false, $true, 1, $false, "41+58=99", 2) | Out-Null
$d.Content.Find.Execute("60-55=5", $true, $false, $false, $false, $false, $true, 1, $false, "2+22=24", 2) | Out-Null
$d.Content.Find.Execute("57-1=56", $true, $false, $false, $false, $false, $true, 1, $false, "89+3=92", 2) | Out-Null
$d.Content.Find.Execute("77-42=35", $true, $false, $false, $false, $false, $true, 1, $false, "0+99=99", 2) | Out-Null
$d.Content.Find.Execute("81+12=93", $true, $false, $false, $false, $false, $true, 1, $false, "67-1=66", 2) | Out-Null
$d.Content.Find.Execute("35-12=23", $true, $false, $false, $false, $false, $true, 1, $false, "26+7=33", 2) | Out-Null
$d.Content.Find.Execute("61-26=35", $true, $false, $false, $false, $false, $true, 1, $false, "96+2=98", 2) | Out-Null
$d.Content.Find.Execute("78-5=73", $true, $false, $false, $false, $false, $true, 1, $false, "16-7=9", 2) | Out-Null
$d.Content.Find.Execute("25-6=19", $true, $false, $false, $false, $false, $true, 1, $false, "67-40=27", 2) | Out-Null
$d.Content.Find.Execute("99-41=58", $true, $false, $false, $false, $false, $true, 1, $false, "69-39=30", 2) | Out-Null
$d.Content.Find.Execute("10+35=45", $true, $false, $false, $false, $false, $true, 1, $false, "89-60=29", 2) | Out-Null
$d.Content.Find.Execute("7+46=53", $true, $false, $false, $false, $false, $true, 1, $false, "64-57=7", 2) | Out-Null
$d.Content.Find.Execute("74+11=85", $true, $false, $false, $false, $false, $true, 1, $false, "3+70=73", 2) | Out-Null
$d.Content.Find.Execute("30+6=36", $true, $false, $false, $false, $false, $true, 1, $false, "46-15=31", 2) | Out-Null
$d.Content.Find.Execute("66-54=12", $true, $false, $false, $false, $false, $true, 1, $false, "17+38=55", 2) | Out-Null
$d.Content.Find.Execute("10+19=29", $true, $false, $false, $false, $false, $true, 1, $false, "51-46=5", 2) | Out-Null
$d.Content.Find.Execute("34+60=94", $true, $false, $false, $false, $false, $true, 1, $false, "98-50=48", 2) | Out-Null
$d.Content.Find.Execute("60+38=98", $true, $false, $false, $false, $false, $true, 1, $false, "77-2=75", 2) | Out-Null
$d.Content.Find.Execute("25+51=76", $true, $false, $false, $false, $false, $true, 1, $false, "75-57=18", 2) | Out-Null
$d.Content.Find.Execute("77-25=52", $true, $false, $false, $false, $false, $true, 1, $false, "31+13=44", 2) | Out-Null
$d.Content.Find.Execute("48-27=21", $true, $false, $false, $false, $false, $true, 1, $false, "44-22=22", 2) | Out-Null
$d.Content.Find.Execute("47-0=47", $true, $false, $false, $false, $false, $true, 1, $false, "19+9=28", 2) | Out-Null
$d.Content.Find.Execute("19+58=77", $true, $false, $false, $false, $false, $true, 1, $false, "8+12=20", 2) | Out-Null
$d.Content.Find.Execute("17-3=14", $true, $false, $false, $false, $false, $true, 1, $false, "95-35=60", 2) | Out-Null
$d.Content.Find.Execute("16+63=79", $true, $false, $false, $false, $false, $true, 1, $false, "7+17=24", 2) | Out-Null
$d.Content.Find.Execute("52-8=44", $true, $false, $false, $false, $false, $true, 1, $false, "22+36=58", 2) | Out-Null
$d.Content.Find.Execute("65-26=39", $true, $false, $false, $false, $false, $true, 1, $false, "5+56=61", 2) | Out-Null
$d.Content.Find.Execute("54-31=23", $true, $false, $false, $false, $false, $true, 1, $false, "35+47=82", 2) | Out-Null
$d.Content.Find.Execute("62-20=42", $true, $false, $false, $false, $false, $true, 1, $false, "86-35=51", 2) | Out-Null
$d.Content.Find.Execute("98-32=66", $true, $false, $false, $false, $false, $true, 1, $false, "85-28=57", 2) | Out-Null
$d.Content.Find.Execute("74-42=32", $true, $false, $false, $false, $false, $true, 1, $false, "73-39=34", 2) | Out-Null
$d.Content.Find.Execute("57-25=32", $true, $false, $false, $false, $false, $true, 1, $false, "93-80=13", 2) | Out-Null
$d.Content.Find.Execute("64+28=92", $true, $false, $false, $false, $false, $true, 1, $false, "54-24=30", 2) | Out-Null
$d.Content.Find.Execute("31-31=0", $true, $false, $false, $false, $false, $true, 1, $false, "60+11=71", 2) | Out-Null
$d.Content.Find.Execute("51+44=95", $true, $false, $false, $false, $false, $true, 1, $false, "44+24=68", 2) | Out-Null
$d.Content.Find.Execute("22-17=5", $true, $false, $false, $false, $false, $true, 1, $false, "8+66=74", 2) | Out-Null
$d.Content.Find.Execute("2+67=69", $true, $false, $false, $false, $false, $true, 1, $false, "60+11=71", 2) | Out-Null
$d.Content.Find.Execute("67-61=6", $true, $false, $false, $false, $false, $true, 1, $false, "92-33=59", 2) | Out-Null
$d.Content.Find.Execute("99-85=14", $true, $false, $false, $false, $false, $true, 1, $false, "14-2=12", 2) | Out-Null
$d.Content.Find.Execute("67-0=67", $true, $false, $false, $false, $false, $true, 1, $false, "76-64=12", 2) | Out-Null
$d.Content.Find.Execute("98-64=34", $true, $false, $false, $false, $false, $true, 1, $false, "16+33=49", 2) | Out-Null
$d.Content.Find.Execute("28-19=9", $true, $false, $false, $false, $false, $true, 1, $false, "26+72=98", 2) | Out-Null
$d.Content.Find.Execute("89-80=9", $true, $false, $false, $false, $false, $true, 1, $false, "58+15=73", 2) | Out-Null
$d.Content.Find.Execute("56-44=12", $true, $false, $false, $false, $false, $true, 1, $false, "53+42=95", 2) | Out-Null
$d.Content.Find.Execute("97-44=53", $true, $false, $false, $false, $false, $true, 1, $false, "82-79=3", 2) | Out-Null
$d.Content.Find.Execute("92+0=92", $true, $false, $false, $false, $false, $true, 1, $false, "30-21=9", 2) | Out-Null
$d.Content.Find.Execute("27-8=19", $true, $false, $false, $false, $false, $true, 1, $false, "86-58=28", 2) | Out-Null
$d.Content.Find.Execute("82-21=61", $true, $false, $false, $false, $false, $true, 1, $false, "70-20=50", 2) | Out-Null
$d.Content.Find.Execute("55-51=4", $true, $false, $false, $false, $false, $true, 1, $false, "24+33=57", 2) | Out-Null
$d.Content.Find.Execute("23+7=30", $true, $false, $false, $false, $false, $true, 1, $false, "94+1=95", 2) | Out-Null
$d.Content.Find.Execute("9+18=27", $true, $false, $false, $false, $false, $true, 1, $false, "3+41=44", 2) | Out-Null
$d.Content.Find.Execute("46-17=29", $true, $false, $false, $false, $false, $true, 1, $false, "9+46=55", 2) | Out-Null
$d.Content.Find.Execute("82-33=49", $true, $false, $false, $false, $false, $true, 1, $false, "71+6=77", 2) | Out-Null
$d.Content.Find.Execute("88-65=23", $true, $false, $false, $false, $false, $true, 1, $false, "59-39=20", 2) | Out-Null
$d.Content.Find.Execute("91-89=2", $true, $false, $false, $false, $false, $true, 1, $false, "71+12=83", 2) | Out-Null
$d.Content.Find.Execute("77-26=51", $true, $false, $false, $false, $false, $true, 1, $false, "76-54=22", 2) | Out-Null
$d.Content.Find.Execute("47-18=29", $true, $false, $false, $false, $false, $true, 1, $false, "47-42=5", 2) | Out-Null
$d.Content.Find.Execute("14+64=78", $true, $false, $false, $false, $false, $true, 1, $false, "38-23=15", 2) | Out-Null
$d.Content.Find.Execute("57+5=62", $true, $false, $false, $false, $false, $true, 1, $false, "50+8=58", 2) | Out-Null
$d.Content.Find.Execute("70-51=19", $true, $false, $false, $false, $false, $true, 1, $false, "98-19=79", 2) | Out-Null
$d.Content.Find.Execute("68+12=80", $true, $false, $false, $false, $false, $true, 1, $false, "22+20=42", 2) | Out-Null
$d.Content.Find.Execute("1+32=33", $true, $false, $false, $false, $false, $true, 1, $false, "57+38=95", 2) | Out-Null
$d.Content.Find.Execute("0+53=53", $true, $false, $false, $false, $false, $true, 1, $false, "85-28=57", 2) | Out-Null
$d.Content.Find.Execute("43+34=77", $true, $false, $false, $false, $false, $true, 1, $false, "29+17=46", 2) | Out-Null
$d.Content.Find.Execute("86-81=5", $true, $false, $false, $false, $false, $true, 1, $false, "93-47=46", 2) | Out-Null
$d.Content.Find.Execute("7+35=42", $true, $false, $false, $false, $false, $true, 1, $false, "84-13=71", 2) | Out-Null
$d.Content.Find.Execute("97-12=85", $true, $false, $false, $false, $false, $true, 1, $false, "49-43=6", 2) | Out-Null
$d.Content.Find.Execute("15+64=79", $true, $false, $false, $false, $false, $true, 1, $false, "92-6=86", 2) | Out-Null
